$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while guaranteeing it stays a text
# value (matching the original inlineStr cells), even when the text looks
# like a plain decimal number (e.g. "1.00", "5.99"). Cells whose text is
# not a number will be unaffected by the NumberFormat dance, and the
# final Style reset keeps the cell's visual style identical to before
# (no explicit style index), just like the original file.
function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Rows 29 and 30 swap places (Monero <-> Toncoin) plus value updates.
$ws.Cells.Item(29, 2).Value = "Toncoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Cells.Item(29, 4) "2.07"
$ws.Cells.Item(29, 5).Value = "  -4.92%  "

$ws.Cells.Item(30, 2).Value = "Monero"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Cells.Item(30, 4) "167.61"
$ws.Cells.Item(30, 5).Value = "  +1.71%  "

# Price (D) and Volume(1h) (E) column updates for all other changed rows.
$updates = @(
    @{ Row = 2;  D = "42.386.14";   E = "  -0.74%  " },
    @{ Row = 3;  D = "2.280.08";    E = "  -1.05%  " },
    @{ Row = 4;  D = "1.00";        E = "  +0.03%  " },
    @{ Row = 5;  D = "300.39";      E = "  -0.68%  " },
    @{ Row = 6;  D = "96.35";       E = "  -2.79%  " },
    @{ Row = 7;  D = "0.497";       E = "  -0.94%  " },
    @{ Row = 9;  D = "0.492";       E = "  -1.70%  " },
    @{ Row = 10; D = "33.33";       E = "  -3.90%  " },
    @{ Row = 11; D = "0.0790";      E = "  -0.02%  " },
    @{ Row = 12; D = "48.25";       E = "  -6.11%  " },
    @{ Row = 13; E = "  +2.12%  " },
    @{ Row = 14; D = "16.00";       E = "  +2.55%  " },
    @{ Row = 15; E = "  +0.13%  " },
    @{ Row = 16; D = "2.634.18";    E = "  -0.88%  " },
    @{ Row = 17; D = "2.288.04";    E = "  -0.44%  " },
    @{ Row = 18; E = "  -0.65%  " },
    @{ Row = 19; D = "42.322.16";   E = "  -0.67%  " },
    @{ Row = 20; E = "  +1.53%  " },
    @{ Row = 21; E = "  -1.01%  " },
    @{ Row = 22; D = "5.99";        E = "  -0.68%  " },
    @{ Row = 23; E = "  -1.93%  " },
    @{ Row = 24; D = "235.63";      E = "  +0.17%  " },
    @{ Row = 25; D = "1.96";        E = "  +0.74%  " },
    @{ Row = 26; E = "  -0.03%  " },
    @{ Row = 27; E = "  -2.05%  " },
    @{ Row = 28; D = "23.79";       E = "  -3.60%  " },
    @{ Row = 31; D = "33.72";       E = "  -1.51%  " },
    @{ Row = 32; E = "  +0.52%  " },
    @{ Row = 33; D = "1.00";        E = "  +0.03%  " },
    @{ Row = 34; D = "4.66";        E = "  +4.88%  " },
    @{ Row = 35; E = "  -1.29%  " },
    @{ Row = 36; D = "16.81";       E = "  +0.83%  " },
    @{ Row = 37; E = "  -3.10%  " },
    @{ Row = 38; D = "0.0693";      E = "  -0.13%  " },
    @{ Row = 39; E = "  -3.30%  " },
    @{ Row = 40; D = "0.0993";      E = "  -0.81%  " },
    @{ Row = 41; E = "  -3.58%  " },
    @{ Row = 42; E = "  -1.70%  " },
    @{ Row = 43; D = "2.25";        E = "  -8.51%  " },
    @{ Row = 44; D = "1.957.88";    E = "  -0.23%  " },
    @{ Row = 45; E = "  +0.12%  " },
    @{ Row = 46; D = "17.62";       E = "  -4.30%  " },
    @{ Row = 47; D = "9.63";        E = "  -5.57%  " },
    @{ Row = 48; E = "  -2.44%  " },
    @{ Row = 49; D = "2.503.84";    E = "  -0.94%  " },
    @{ Row = 50; D = "52.44";       E = "  -5.24%  " },
    @{ Row = 51; E = "  -3.02%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
